# User Story BL054: add "Take 3!" results block (rows 32-41) to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 32: section label ---
$ws.Range("A32").Value = "Take 3!"

# --- Row 33: header row ---
$ws.Range("B33").Value = "MeanDeltaTime"
$ws.Range("C33").Value = "CorrectedTime"
$ws.Range("D33").Value = "DirectDistance"
$ws.Range("E33").Value = "CalculatedDistance"
$ws.Range("F33").Value = "PWV (m/s)"
$ws.Range("G33").Value = "StandardDeviation"
$ws.Range("H33").Value = "IsStandardDeviationValid"
$ws.Range("I33").Value = "HR"
$ws.Range("J33").Value = "SignalLength"
$ws.Range("K33").Value = "Deltas"
$ws.Range("L33").Value = "Valid Deltas"

# --- Row 34: PWV_IK2 raw values ---
$ws.Range("A34").Value = "PWV_IK2"
$ws.Range("B34").Value = 82
$ws.Range("C34").Value = 42
$ws.Range("D34").Value = 600
$ws.Range("E34").Value = 400
$ws.Range("F34").Value = 9.6
$ws.Range("G34").Value = 0.4
$ws.Range("H34").Value = $true
$ws.Range("I34").Value = 56
$ws.Range("J34").Value = 2730
$ws.Range("K34").Value = "?"

# --- Row 35: PWV_DL1 raw values ---
$ws.Range("A35").Value = "PWV_DL1"
$ws.Range("B35").Value = 68
$ws.Range("C35").Value = 28
$ws.Range("D35").Value = 530
$ws.Range("E35").Value = 330
$ws.Range("F35").Value = 11.7
$ws.Range("G35").Value = 0.5
$ws.Range("H35").Value = $true
$ws.Range("I35").Value = 65
$ws.Range("J35").Value = 2350
$ws.Range("K35").Value = "?"

# --- Row 36: PWV_09 raw values ---
$ws.Range("A36").Value = "PWV_09"
$ws.Range("B36").Value = 126
$ws.Range("C36").Value = 87
$ws.Range("D36").Value = 740
$ws.Range("E36").Value = 540
$ws.Range("F36").Value = 6.3
$ws.Range("G36").Value = 0.7
$ws.Range("H36").Value = $false
$ws.Range("I36").Value = 58
$ws.Range("J36").Value = 2640
$ws.Range("K36").Value = "?"

# --- Row 38: Results label ---
$ws.Range("A38").Value = "Results"

# --- Row 39: PWV_IK2 computed results ---
$ws.Range("A39").Value = "PWV_IK2"
$ws.Range("B39").Value = 80.861557000000005
$ws.Range("C39").Value = 40.861556999999998
$ws.Range("D39").Value = 600
$ws.Range("E39").Value = 400
$ws.Range("F39").Value = 9.8067969999999995
$ws.Range("G39").Value = 0.42628961999999998
$ws.Range("H39").Value = $true
$ws.Range("I39").Value = 56.263924000000003
$ws.Range("J39").Value = 2730
$ws.Range("K39").Value = 8
$ws.Range("L39").Value = 8

# --- Row 40: PWV_DL1 computed results ---
$ws.Range("A40").Value = "PWV_DL1"
$ws.Range("B40").Value = 67.484497000000005
$ws.Range("C40").Value = 27.484497000000001
$ws.Range("D40").Value = 530
$ws.Range("E40").Value = 330
$ws.Range("F40").Value = 12.061915000000001
$ws.Range("G40").Value = 0.80147111000000004
$ws.Range("H40").Value = $true
$ws.Range("I40").Value = 65.327667000000005
$ws.Range("J40").Value = 2350
$ws.Range("K40").Value = 8
$ws.Range("L40").Value = 8

# --- Row 41: PWV_09 computed results ---
$ws.Range("A41").Value = "PWV_09"
$ws.Range("B41").Value = 124.31859
$ws.Range("C41").Value = 84.318588000000005
$ws.Range("D41").Value = 740
$ws.Range("E41").Value = 540
$ws.Range("F41").Value = 6.4872775000000003
$ws.Range("G41").Value = 0.71557104999999999
$ws.Range("H41").Value = $false
$ws.Range("I41").Value = 58.082222000000002
$ws.Range("J41").Value = 2640
$ws.Range("K41").Value = 8
$ws.Range("L41").Value = 8

# --- Update selection to match the new view state ---
$ws.Range("G42").Select()
